# Update automatic: dades i banners [2026-02-10 20:50]
# Applies the per-cell value updates described by the diff against resum_diari_meteocat.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (outside the A1:O46 data range) used to stage pure-percentage text values
# so Excel keeps them as literal text (e.g. "97%") instead of auto-converting them to a
# numeric percentage (0.97) with a new number format / style.
$helper = $ws.Range("Q1")

# E2: '2026-02-10 20:18:25' -> '2026-02-10 20:48:26'
$ws.Range("E2").Value = '2026-02-10 20:48:26'
# I2: '46.7 mm' -> '46.8 mm'
$ws.Range("I2").Value = '46.8 mm'
# E3: '2026-02-10 20:18:27' -> '2026-02-10 20:48:28'
$ws.Range("E3").Value = '2026-02-10 20:48:28'
# H3: '96%' -> '97%'
$helper.NumberFormat = "@"
$helper.Value = '97%'
$helper.Copy()
$ws.Range("H3").PasteSpecial(-4163)
$helper.Clear()
# E4: '2026-02-10 20:18:30' -> '2026-02-10 20:48:31'
$ws.Range("E4").Value = '2026-02-10 20:48:31'
# E5: '2026-02-10 20:18:32' -> '2026-02-10 20:48:33'
$ws.Range("E5").Value = '2026-02-10 20:48:33'
# E6: '2026-02-10 20:18:34' -> '2026-02-10 20:48:36'
$ws.Range("E6").Value = '2026-02-10 20:48:36'
# I6: '0.2 mm' -> '0.9 mm'
$ws.Range("I6").Value = '0.9 mm'
# L6: '13.3 km/h - 347º 19:52 TU' -> '22.7 km/h - 333º 20:22 TU'
$ws.Range("L6").Value = '22.7 km/h - 333º 20:22 TU'
# E7: '2026-02-10 20:18:37' -> '2026-02-10 20:48:38'
$ws.Range("E7").Value = '2026-02-10 20:48:38'
# E8: '2026-02-10 20:18:39' -> '2026-02-10 20:48:40'
$ws.Range("E8").Value = '2026-02-10 20:48:40'
# J8: '1004.4 hPa' -> '1004.3 hPa'
$ws.Range("J8").Value = '1004.3 hPa'
# O8: '12.1 °C' -> '12.2 °C'
$ws.Range("O8").Value = '12.2 °C'
# E9: '2026-02-10 20:18:42' -> '2026-02-10 20:48:43'
$ws.Range("E9").Value = '2026-02-10 20:48:43'
# I9: '0.7 mm' -> '3.1 mm'
$ws.Range("I9").Value = '3.1 mm'
# E10: '2026-02-10 20:18:44' -> '2026-02-10 20:48:45'
$ws.Range("E10").Value = '2026-02-10 20:48:45'
# O10: '10.2 °C' -> '10.3 °C'
$ws.Range("O10").Value = '10.3 °C'
# E11: '2026-02-10 20:18:47' -> '2026-02-10 20:48:48'
$ws.Range("E11").Value = '2026-02-10 20:48:48'
# O11: '7.7 °C' -> '7.8 °C'
$ws.Range("O11").Value = '7.8 °C'
# E12: '2026-02-10 20:18:49' -> '2026-02-10 20:48:50'
$ws.Range("E12").Value = '2026-02-10 20:48:50'
# I12: '0.6 mm' -> '1.7 mm'
$ws.Range("I12").Value = '1.7 mm'
# E13: '2026-02-10 20:18:51' -> '2026-02-10 20:48:52'
$ws.Range("E13").Value = '2026-02-10 20:48:52'
# H13: '90%' -> '91%'
$helper.NumberFormat = "@"
$helper.Value = '91%'
$helper.Copy()
$ws.Range("H13").PasteSpecial(-4163)
$helper.Clear()
# I13: '9.3 mm' -> '10.5 mm'
$ws.Range("I13").Value = '10.5 mm'
# J13: '1006.7 hPa' -> '1006.6 hPa'
$ws.Range("J13").Value = '1006.6 hPa'
# O13: '5.2 °C' -> '5.3 °C'
$ws.Range("O13").Value = '5.3 °C'
# E14: '2026-02-10 20:18:54' -> '2026-02-10 20:48:55'
$ws.Range("E14").Value = '2026-02-10 20:48:55'
# H14: '86%' -> '87%'
$helper.NumberFormat = "@"
$helper.Value = '87%'
$helper.Copy()
$ws.Range("H14").PasteSpecial(-4163)
$helper.Clear()
# E15: '2026-02-10 20:18:56' -> '2026-02-10 20:48:57'
$ws.Range("E15").Value = '2026-02-10 20:48:57'
# I15: '0.9 mm' -> '3.5 mm'
$ws.Range("I15").Value = '3.5 mm'
# O15: '9.0 °C' -> '9.1 °C'
$ws.Range("O15").Value = '9.1 °C'
# E16: '2026-02-10 20:18:59' -> '2026-02-10 20:49:00'
$ws.Range("E16").Value = '2026-02-10 20:49:00'
# I16: '25.5 mm' -> '25.9 mm'
$ws.Range("I16").Value = '25.9 mm'
# E17: '2026-02-10 20:19:01' -> '2026-02-10 20:49:02'
$ws.Range("E17").Value = '2026-02-10 20:49:02'
# I17: '0.1 mm' -> '0.2 mm'
$ws.Range("I17").Value = '0.2 mm'
# M17: '6.6 °C 10:07 TU' -> '6.7 °C 20:06 TU'
$ws.Range("M17").Value = '6.7 °C 20:06 TU'
# O17: '4.7 °C' -> '4.8 °C'
$ws.Range("O17").Value = '4.8 °C'
# E18: '2026-02-10 20:19:03' -> '2026-02-10 20:49:05'
$ws.Range("E18").Value = '2026-02-10 20:49:05'
# J18: '1004.3 hPa' -> '1004.2 hPa'
$ws.Range("J18").Value = '1004.2 hPa'
# E19: '2026-02-10 20:19:06' -> '2026-02-10 20:49:07'
$ws.Range("E19").Value = '2026-02-10 20:49:07'
# I19: '0.1 mm' -> '0.4 mm'
$ws.Range("I19").Value = '0.4 mm'
# O19: '6.7 °C' -> '6.8 °C'
$ws.Range("O19").Value = '6.8 °C'
# E20: '2026-02-10 20:19:08' -> '2026-02-10 20:49:09'
$ws.Range("E20").Value = '2026-02-10 20:49:09'
# I20: '8.4 mm' -> '9.8 mm'
$ws.Range("I20").Value = '9.8 mm'
# E21: '2026-02-10 20:19:10' -> '2026-02-10 20:49:12'
$ws.Range("E21").Value = '2026-02-10 20:49:12'
# I21: '8.7 mm' -> '8.8 mm'
$ws.Range("I21").Value = '8.8 mm'
# E22: '2026-02-10 20:19:13' -> '2026-02-10 20:49:14'
$ws.Range("E22").Value = '2026-02-10 20:49:14'
# E23: '2026-02-10 20:19:15' -> '2026-02-10 20:49:17'
$ws.Range("E23").Value = '2026-02-10 20:49:17'
# E24: '2026-02-10 20:19:18' -> '2026-02-10 20:49:19'
$ws.Range("E24").Value = '2026-02-10 20:49:19'
# J24: '1006.0 hPa' -> '1005.9 hPa'
$ws.Range("J24").Value = '1005.9 hPa'
# E25: '2026-02-10 20:19:20' -> '2026-02-10 20:49:21'
$ws.Range("E25").Value = '2026-02-10 20:49:21'
# I25: '19.6 mm' -> '22.2 mm'
$ws.Range("I25").Value = '22.2 mm'
# E26: '2026-02-10 20:19:23' -> '2026-02-10 20:49:24'
$ws.Range("E26").Value = '2026-02-10 20:49:24'
# E27: '2026-02-10 20:19:25' -> '2026-02-10 20:49:26'
$ws.Range("E27").Value = '2026-02-10 20:49:26'
# I27: '11.9 mm' -> '12.1 mm'
$ws.Range("I27").Value = '12.1 mm'
# E28: '2026-02-10 20:19:28' -> '2026-02-10 20:49:29'
$ws.Range("E28").Value = '2026-02-10 20:49:29'
# I28: '1.3 mm' -> '1.9 mm'
$ws.Range("I28").Value = '1.9 mm'
# E29: '2026-02-10 20:19:30' -> '2026-02-10 20:49:31'
$ws.Range("E29").Value = '2026-02-10 20:49:31'
# E30: '2026-02-10 20:19:33' -> '2026-02-10 20:49:34'
$ws.Range("E30").Value = '2026-02-10 20:49:34'
# I30: '0.5 mm' -> '1.4 mm'
$ws.Range("I30").Value = '1.4 mm'
# J30: '1004.3 hPa' -> '1004.2 hPa'
$ws.Range("J30").Value = '1004.2 hPa'
# E31: '2026-02-10 20:19:35' -> '2026-02-10 20:49:36'
$ws.Range("E31").Value = '2026-02-10 20:49:36'
# I31: '1.2 mm' -> '3.2 mm'
$ws.Range("I31").Value = '3.2 mm'
# E32: '2026-02-10 20:19:38' -> '2026-02-10 20:49:39'
$ws.Range("E32").Value = '2026-02-10 20:49:39'
# O32: '10.5 °C' -> '10.6 °C'
$ws.Range("O32").Value = '10.6 °C'
# E33: '2026-02-10 20:19:40' -> '2026-02-10 20:49:41'
$ws.Range("E33").Value = '2026-02-10 20:49:41'
# I33: '11.0 mm' -> '12.0 mm'
$ws.Range("I33").Value = '12.0 mm'
# E34: '2026-02-10 20:19:43' -> '2026-02-10 20:49:43'
$ws.Range("E34").Value = '2026-02-10 20:49:43'
# I34: '14.2 mm' -> '15.2 mm'
$ws.Range("I34").Value = '15.2 mm'
# E35: '2026-02-10 20:19:45' -> '2026-02-10 20:49:46'
$ws.Range("E35").Value = '2026-02-10 20:49:46'
# E36: '2026-02-10 20:19:47' -> '2026-02-10 20:49:48'
$ws.Range("E36").Value = '2026-02-10 20:49:48'
# I36: '0.9 mm' -> '2.7 mm'
$ws.Range("I36").Value = '2.7 mm'
# E37: '2026-02-10 20:19:50' -> '2026-02-10 20:49:51'
$ws.Range("E37").Value = '2026-02-10 20:49:51'
# I37: '3.3 mm' -> '3.5 mm'
$ws.Range("I37").Value = '3.5 mm'
# E38: '2026-02-10 20:19:52' -> '2026-02-10 20:49:53'
$ws.Range("E38").Value = '2026-02-10 20:49:53'
# H38: '88%' -> '87%'
$helper.NumberFormat = "@"
$helper.Value = '87%'
$helper.Copy()
$ws.Range("H38").PasteSpecial(-4163)
$helper.Clear()
# E39: '2026-02-10 20:19:55' -> '2026-02-10 20:49:56'
$ws.Range("E39").Value = '2026-02-10 20:49:56'
# I39: '10.6 mm' -> '11.2 mm'
$ws.Range("I39").Value = '11.2 mm'
# E40: '2026-02-10 20:19:57' -> '2026-02-10 20:49:58'
$ws.Range("E40").Value = '2026-02-10 20:49:58'
# I40: '13.4 mm' -> '14.4 mm'
$ws.Range("I40").Value = '14.4 mm'
# J40: '1007.0 hPa' -> '1006.9 hPa'
$ws.Range("J40").Value = '1006.9 hPa'
# O40: '7.5 °C' -> '7.6 °C'
$ws.Range("O40").Value = '7.6 °C'
# E41: '2026-02-10 20:20:00' -> '2026-02-10 20:50:00'
$ws.Range("E41").Value = '2026-02-10 20:50:00'
# H41: '82%' -> '81%'
$helper.NumberFormat = "@"
$helper.Value = '81%'
$helper.Copy()
$ws.Range("H41").PasteSpecial(-4163)
$helper.Clear()
# M41: '20.2 °C 10:34 TU' -> '21.0 °C 20:28 TU'
$ws.Range("M41").Value = '21.0 °C 20:28 TU'
# O41: '14.3 °C' -> '14.4 °C'
$ws.Range("O41").Value = '14.4 °C'
# E42: '2026-02-10 20:20:02' -> '2026-02-10 20:50:03'
$ws.Range("E42").Value = '2026-02-10 20:50:03'
# E43: '2026-02-10 20:20:04' -> '2026-02-10 20:50:05'
$ws.Range("E43").Value = '2026-02-10 20:50:05'
# O43: '9.7 °C' -> '9.8 °C'
$ws.Range("O43").Value = '9.8 °C'
# E44: '2026-02-10 20:20:07' -> '2026-02-10 20:50:08'
$ws.Range("E44").Value = '2026-02-10 20:50:08'
# I44: '28.3 mm' -> '28.8 mm'
$ws.Range("I44").Value = '28.8 mm'
# E45: '2026-02-10 20:20:09' -> '2026-02-10 20:50:10'
$ws.Range("E45").Value = '2026-02-10 20:50:10'
# E46: '2026-02-10 20:20:11' -> '2026-02-10 20:50:12'
$ws.Range("E46").Value = '2026-02-10 20:50:12'
# H46: '81%' -> '80%'
$helper.NumberFormat = "@"
$helper.Value = '80%'
$helper.Copy()
$ws.Range("H46").PasteSpecial(-4163)
$helper.Clear()
# J46: '1005.9 hPa' -> '1005.8 hPa'
$ws.Range("J46").Value = '1005.8 hPa'
# O46: '14.5 °C' -> '14.7 °C'
$ws.Range("O46").Value = '14.7 °C'
